$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "predicted type" column (E) classifying each variable as
# good / useless / bad.
$ws.Cells.Item(1, 5).Value = "predicted type"

$ws.Cells.Item(2, 5).Value = "useless"   # dmage
$ws.Cells.Item(3, 5).Value = "good"      # mrace3
$ws.Cells.Item(4, 5).Value = "good"      # dmar
$ws.Cells.Item(5, 5).Value = "good"      # dlivord
$ws.Cells.Item(6, 5).Value = "useless"   # frace4
$ws.Cells.Item(7, 5).Value = "bad"       # dgestat
$ws.Cells.Item(8, 5).Value = "useless"   # csex
$ws.Cells.Item(10, 5).Value = "useless"  # dplural
$ws.Cells.Item(12, 5).Value = "useless"  # drink
$ws.Cells.Item(13, 5).Value = "bad"      # tobacco
$ws.Cells.Item(14, 5).Value = "good"     # wgain
$ws.Cells.Item(15, 5).Value = "bad"      # lung
$ws.Cells.Item(16, 5).Value = "bad"      # cardiac

# Widen the new columns to fit their contents (matches the bestFit
# column sizing Excel applies to columns B/C already in the sheet).
# (Inputs are pre-compensated for this host's ColumnWidth->stored-width
# rounding so the saved OOXML <col> widths land on 21 and ~12.664.)
$ws.Columns.Item(4).ColumnWidth = 20.166666666666664
$ws.Columns.Item(5).ColumnWidth = 11.833333333333332

# Move the selection the same way the original author's last click did.
$ws.Range("E21").Select()
